# clientes.xlsx — "Atualizar código com melhorias"
#
# The "vencimento" column (C) used to hold dates typed as plain text
# ("30/8/2025", "20/8/2025", "10/8/2025") under a text ("@") number
# format. This converts those three cells to real Excel date serials
# displayed as dd/mm/yyyy, and leaves the selection on E6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("clientes")

# Switch column C's (rows 2-4) number format from text to a date format
# BEFORE writing the values, so the dates are stored/recognised as real
# numeric dates instead of being coerced to text.
$ws.Range("C2:C4").NumberFormat = "dd/mm/yyyy"

$ws.Range("C2").Value = (Get-Date -Year 2025 -Month 8 -Day 30).Date
$ws.Range("C3").Value = (Get-Date -Year 2025 -Month 8 -Day 20).Date
$ws.Range("C4").Value = (Get-Date -Year 2025 -Month 8 -Day 10).Date

# Match the saved selection (E6).
$ws.Range("E6").Select()
